$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Get User Account with wrong Auth") - ResponseMapKeys / ResponseMapValues
# used to hold the "team" mapping; now holds the simple account mapping.
$ws.Range("F3").Value = "account.email_address"
$ws.Range("G3").Value = "aashish.kumar@sofbang.com"

# Row 5 ("Get Team") - ResponseMapKeys / ResponseMapValues
# used to hold the simple account mapping; now holds the expanded team mapping.
$ws.Range("G5").Value = "Sofbang Team;vivek.ahuja@sofbang.com;puneet.gandhi@sofbang.com;5000"
$ws.Range("F5").Value = "team.name;team.accounts[0].email_address;team.accounts[1].email_address;team.accounts[1].quotas.api_signature_requests_left"

# Widen columns F and G to fit the new (longer) content
$ws.Columns.Item(6).ColumnWidth = 78.66666666666667
$ws.Columns.Item(7).ColumnWidth = 37.5

# Update the active selection to reflect where the editor left off
$ws.Range("F9").Select()
